# "Hjemme passive tweaks lichtwark deleted values"
# Update the first four data columns (B:E) on rows 1-3 of Ark1 with the
# corrected / re-measured values, and shrink the lingering cell selection
# from the old full-range (B1:AY3) down to just the edited block (B1:E3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - angle/header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - "CON" subject data
$ws.Range("B2").Value = 5.0637421277820867
$ws.Range("C2").Value = 4.2841422700928815
$ws.Range("D2").Value = 6.9414908877550401
$ws.Range("E2").Value = 6.1128034389697872

# Row 3 - "STR" subject data
$ws.Range("B3").Value = 4.4550762181419969
$ws.Range("C3").Value = 6.8061620425162186
$ws.Range("D3").Value = 7.2646165724020548
$ws.Range("E3").Value = 5.5698631668856535

# Leave the selection on the touched block, matching the saved view state.
$ws.Range("B1:E3").Select() | Out-Null
